# "added ne scripts for ENWIAM" -- append two new test-case rows
# (ENWIAM00015, ENWIAM00016) to the "Test Cases" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Row 48 : ENWIAM00015 ---------------------------------------------
$ws.Cells.Item(48, 1).Value = "ENWIAM00015"
$ws.Cells.Item(48, 2).Value = "OPQA-1870||OPQA-1874"
$ws.Cells.Item(48, 3).Value = 'Verify that If the STeAM account that is trying to be linked/merged by the user is in a "locked"/Suspended status, then the link/merge shall not be made and the user shall be informed that the STeAM account is locked.'
$ws.Cells.Item(48, 4).Value = "Y"

# --- Row 49 : ENWIAM00016 ---------------------------------------------
$ws.Cells.Item(49, 1).Value = "ENWIAM00016"
$ws.Cells.Item(49, 2).Value = "OPQA-2362||OPQA-2359"
$ws.Cells.Item(49, 3).Value = "Verify that upon a successful sign-in for the first time on the ENW landing screen using STeAM, the user shall be prompted to link existing Neon accounts that have the same email address as the newly registered account"
$ws.Cells.Item(49, 4).Value = "Y"

# New font/style for the Description column of these two rows: dark-grey
# (#333333) wrapped text, no cell border (distinct from the bordered style
# used elsewhere in the table). Apply to C48 first, then clone the
# resulting style onto C49 via a format-only paste so both cells land on
# the very same style record.
$c48 = $ws.Cells.Item(48, 3)
$c48.Borders.LineStyle = -4142
$c48.WrapText = $true
$c48.Font.Color = 3355443

$c48.Copy()
$c49 = $ws.Cells.Item(49, 3)
$c49.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row heights for the new rows.
$ws.Rows.Item(48).RowHeight = 45
$ws.Rows.Item(49).RowHeight = 45

# Reselect row 48 (whole row) to match the saved selection/active-cell.
$ws.Activate() | Out-Null
$ws.Rows.Item(48).Select() | Out-Null

Write-Host "Added ENWIAM00015 and ENWIAM00016 test rows."
